$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A38").Value = "Note:"
$ws.Range("B38").Value = "If there is an assumption, the energy per capita of the compensatory country was multiplied with population of the missing country"

$ws.Range("B39").Select()
